# Remove the Transition_Name_Annot and ISTD_Annot sheets, keeping only
# Sample_Annot, to match the new version of MSTemplate_Creator.
$wb = $excel.ActiveWorkbook

$wsTransition = $wb.Worksheets.Item("Transition_Name_Annot")
$wsTransition.Delete()

$wsIstd = $wb.Worksheets.Item("ISTD_Annot")
$wsIstd.Delete()

$wsSample = $wb.Worksheets.Item("Sample_Annot")
$wsSample.Activate()
